$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $Addr, $Text)
    $rng = $Sheet.Range($Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue $ws "D2" ('26.162.40')
Set-TextValue $ws "E2" ('  -4.52%  ')

# Row 3
Set-TextValue $ws "D3" ('1.656.64')
Set-TextValue $ws "E3" ('  -3.21%  ')

# Row 4
Set-TextValue $ws "D4" ('1.007')
Set-TextValue $ws "E4" ('  +0.31%  ')

# Row 5
Set-TextValue $ws "D5" ('218.04')
Set-TextValue $ws "E5" ('  -2.84%  ')

# Row 6
Set-TextValue $ws "D6" ('0.5161')
Set-TextValue $ws "E6" ('  -3.21%  ')

# Row 7
Set-TextValue $ws "D7" ('1.009')
Set-TextValue $ws "E7" ('  +0.47%  ')

# Row 8
Set-TextValue $ws "D8" ('0.06423')
Set-TextValue $ws "E8" ('  -2.83%  ')

# Row 9
Set-TextValue $ws "D9" ('0.2565')
Set-TextValue $ws "E9" ('  -3.89%  ')

# Row 10
Set-TextValue $ws "D10" ('19.87')
Set-TextValue $ws "E10" ('  -5.04%  ')

# Row 11
Set-TextValue $ws "D11" ('0.07753')
Set-TextValue $ws "E11" ('  +1.41%  ')

# Row 12
Set-TextValue $ws "B12" ('Polkadot')
Set-TextValue $ws "C12" ('https://coinranking.com/coin/25W7FG7om+polkadot-dot')
Set-TextValue $ws "D12" ('4.303')
Set-TextValue $ws "E12" ('  -5.77%  ')

# Row 13
Set-TextValue $ws "D13" ('1.885.86')
Set-TextValue $ws "E13" ('  -3.20%  ')

# Row 14
Set-TextValue $ws "B14" ('WrappedEther')
Set-TextValue $ws "C14" ('https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth')
Set-TextValue $ws "D14" ('1.648.39')
Set-TextValue $ws "E14" ('  -3.87%  ')

# Row 15
Set-TextValue $ws "D15" ('0.5540')
Set-TextValue $ws "E15" ('  -4.08%  ')

# Row 16
Set-TextValue $ws "D16" ('0.0' + ([char]0x2085).ToString() + '8024')
Set-TextValue $ws "E16" ('  -2.12%  ')

# Row 17
Set-TextValue $ws "D17" ('64.39')
Set-TextValue $ws "E17" ('  -5.21%  ')

# Row 18
Set-TextValue $ws "D18" ('26.204.03')
Set-TextValue $ws "E18" ('  -4.40%  ')

# Row 19
Set-TextValue $ws "D19" ('1.012')
Set-TextValue $ws "E19" ('  +0.76%  ')

# Row 20
Set-TextValue $ws "D20" ('210.15')
Set-TextValue $ws "E20" ('  -3.10%  ')

# Row 21
Set-TextValue $ws "D21" ('4.401')
Set-TextValue $ws "E21" ('  -5.88%  ')

# Row 22
Set-TextValue $ws "D22" ('10.09')
Set-TextValue $ws "E22" ('  -3.73%  ')

# Row 23
Set-TextValue $ws "D23" ('5.872')
Set-TextValue $ws "E23" ('  -1.89%  ')

# Row 24
Set-TextValue $ws "D24" ('1.009')
Set-TextValue $ws "E24" ('  +0.46%  ')

# Row 25
Set-TextValue $ws "D25" ('144.44')
Set-TextValue $ws "E25" ('  +1.29%  ')

# Row 26
Set-TextValue $ws "D26" ('1.755')
Set-TextValue $ws "E26" ('  +1.21%  ')

# Row 27
Set-TextValue $ws "D27" ('0.1159')
Set-TextValue $ws "E27" ('  -4.81%  ')

# Row 28
Set-TextValue $ws "D28" ('6.959')
Set-TextValue $ws "E28" ('  -4.56%  ')

# Row 29
Set-TextValue $ws "D29" ('15.78')
Set-TextValue $ws "E29" ('  -3.21%  ')

# Row 30
Set-TextValue $ws "D30" ('0.05253')
Set-TextValue $ws "E30" ('  -3.10%  ')

# Row 31
Set-TextValue $ws "D31" ('1.258')
Set-TextValue $ws "E31" ('  -2.56%  ')

# Row 32
Set-TextValue $ws "D32" ('3.365')
Set-TextValue $ws "E32" ('  -4.23%  ')

# Row 33
Set-TextValue $ws "D33" ('3.215')

# Row 34
Set-TextValue $ws "D34" ('1.567')
Set-TextValue $ws "E34" ('  -4.98%  ')

# Row 35
Set-TextValue $ws "B35" ('MXToken')
Set-TextValue $ws "C35" ('https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx')
Set-TextValue $ws "D35" ('2.748')
Set-TextValue $ws "E35" ('  -4.61%  ')

# Row 36
Set-TextValue $ws "B36" ('HuobiToken')
Set-TextValue $ws "C36" ('https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht')
Set-TextValue $ws "D36" ('2.377')
Set-TextValue $ws "E36" ('  -1.61%  ')

# Row 37
Set-TextValue $ws "D37" ('0.9241')
Set-TextValue $ws "E37" ('  -2.88%  ')

# Row 38
Set-TextValue $ws "D38" ('0.5708')
Set-TextValue $ws "E38" ('  -2.87%  ')

# Row 39
Set-TextValue $ws "D39" ('1.154.29')
Set-TextValue $ws "E39" ('  +10.21%  ')

# Row 40
Set-TextValue $ws "D40" ('0.01592')
Set-TextValue $ws "E40" ('  -2.71%  ')

# Row 41
Set-TextValue $ws "D41" ('1.010')
Set-TextValue $ws "E41" ('  +0.56%  ')

# Row 42
Set-TextValue $ws "D42" ('0.8436')
Set-TextValue $ws "E42" ('  +0.17%  ')

# Row 43
Set-TextValue $ws "E43" ('  -3.59%  ')

# Row 44
Set-TextValue $ws "D44" ('99.91')
Set-TextValue $ws "E44" ('  -1.05%  ')

# Row 45
Set-TextValue $ws "D45" ('1.794.81')
Set-TextValue $ws "E45" ('  -3.26%  ')

# Row 46
Set-TextValue $ws "D46" ('0.0' + ([char]0x2088).ToString() + '109')
Set-TextValue $ws "E46" ('  -3.73%  ')

# Row 47
Set-TextValue $ws "D47" ('0.4499')
Set-TextValue $ws "E47" ('  -0.34%  ')

# Row 48
Set-TextValue $ws "D48" ('55.96')
Set-TextValue $ws "E48" ('  -3.69%  ')

# Row 49
Set-TextValue $ws "D49" ('1.009')
Set-TextValue $ws "E49" ('  +0.50%  ')

# Row 50
Set-TextValue $ws "D50" ('7.897')
Set-TextValue $ws "E50" ('  -2.82%  ')

# Row 51
Set-TextValue $ws "D51" ('0.05092')
Set-TextValue $ws "E51" ('  -2.84%  ')
